$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing row 309 values (low, close, volume)
$ws.Range("E309").Value = 336.916
$ws.Range("F309").Value = 338.665
$ws.Range("G309").Value = 1477207

# Append new rows 310-312 with the latest OANDA:USDHUF OHLCV data
$ws.Range("A310").Value = 45047.29166666666
$ws.Range("B310").Value = "OANDA:USDHUF"
$ws.Range("C310").Value = 338.665
$ws.Range("D310").Value = 351.669
$ws.Range("E310").Value = 335.752
$ws.Range("F310").Value = 346.724
$ws.Range("G310").Value = 1849087

$ws.Range("A311").Value = 45078.29166666666
$ws.Range("B311").Value = "OANDA:USDHUF"
$ws.Range("C311").Value = 346.724
$ws.Range("D311").Value = 347.768
$ws.Range("E311").Value = 335.562
$ws.Range("F311").Value = 341.86
$ws.Range("G311").Value = 1620032

$ws.Range("A312").Value = 45110.29166666666
$ws.Range("B312").Value = "OANDA:USDHUF"
$ws.Range("C312").Value = 341.86
$ws.Range("D312").Value = 357.545
$ws.Range("E312").Value = 340.9
$ws.Range("F312").Value = 350.66
$ws.Range("G312").Value = 380105

# Match the date-time number format/style used in column A for the prior data rows
$ws.Range("A309").Copy()
$ws.Range("A310:A312").PasteSpecial(-4122)
$excel.CutCopyMode = 0
